# Weekly update: insert a new "Choclo" price record for
# "Vega Central Mapocho de Santiago" as row 391, pushing the existing
# rows 391-420 down to 392-421 (dimension grows from R420 to R421).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 391 (formatting is inherited from the row
# above, which is what gives the new D391 cell its date style).
$ws.Rows.Item(391).Insert()

$ws.Range("A391").Value = 9
$ws.Range("B391").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C391").Value = "Metropolitana"
$ws.Range("D391").Value = 44585
$ws.Range("E391").Value = 13
$ws.Range("F391").Value = 100112024
$ws.Range("G391").Value = "Choclo"
$ws.Range("H391").Value = "Choclero"
$ws.Range("I391").Value = "Primera"
$ws.Range("J391").Value = 4300
$ws.Range("K391").Value = 200
$ws.Range("L391").Value = 250
$ws.Range("M391").Value = 225
$ws.Range("N391").Value = "`$/unidad"
$ws.Range("O391").Value = "Región Metropolitana"
$ws.Range("P391").Value = 225
$ws.Range("Q391").Value = 1
$ws.Range("R391").Value = "Hortaliza"
